# Word COM-interop script: fixes mis-inferred PERSON_* placeholder numbers
# in the numbered list (commit: "Oprava inference jmen - castecne reseni
# zkracovani a vlozneho e"). Each numbered list item of the form
#   [[PERSON_x]] - "s [[PERSON_y]]", "o [[PERSON_z]]"
# is corrected to use consistent / renumbered placeholder references.
# The documents paragraph count does not change - this is purely 13
# in-place text replacements, performed in document order so that no
# freshly written text is accidentally re-matched by a later step.

$d = $word.ActiveDocument
$failures = 0

$old1 = ('[[PERSON_93]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_94]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_94]]' + [char]8220)
$new1 = ('[[PERSON_93]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_93]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_93]]' + [char]8220)
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) { $failures = $failures + 1 }
Write-Output "Replace 1: $found1"

$old2 = ('[[PERSON_95]] ' + [char]8211 + ' ' + [char]8222 + 'o [[PERSON_96]]' + [char]8220 + ', ' + [char]8222 + 's [[PERSON_95]]' + [char]8220)
$new2 = ('[[PERSON_94]] ' + [char]8211 + ' ' + [char]8222 + 'o [[PERSON_95]]' + [char]8220 + ', ' + [char]8222 + 's [[PERSON_94]]' + [char]8220)
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $found2) { $failures = $failures + 1 }
Write-Output "Replace 2: $found2"

$old3 = ('[[PERSON_97]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_97]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_97]]' + [char]8220)
$new3 = ('[[PERSON_96]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_96]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_96]]' + [char]8220)
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
if (-not $found3) { $failures = $failures + 1 }
Write-Output "Replace 3: $found3"

$old4 = ('[[PERSON_98]] ' + [char]8211 + ' ' + [char]8222 + 'o [[PERSON_99]]' + [char]8220 + ', ' + [char]8222 + 's [[PERSON_99]]' + [char]8220)
$new4 = ('[[PERSON_97]] ' + [char]8211 + ' ' + [char]8222 + 'o [[PERSON_97]]' + [char]8220 + ', ' + [char]8222 + 's [[PERSON_97]]' + [char]8220)
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
if (-not $found4) { $failures = $failures + 1 }
Write-Output "Replace 4: $found4"

$old5 = ('[[PERSON_100]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_100]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_100]]' + [char]8220)
$new5 = ('[[PERSON_98]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_98]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_98]]' + [char]8220)
$found5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
if (-not $found5) { $failures = $failures + 1 }
Write-Output "Replace 5: $found5"

$old6 = ('[[PERSON_101]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_101]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_102]]' + [char]8220)
$new6 = ('[[PERSON_99]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_99]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_100]]' + [char]8220)
$found6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
if (-not $found6) { $failures = $failures + 1 }
Write-Output "Replace 6: $found6"

$old7 = ('[[PERSON_103]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_103]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_103]]' + [char]8220)
$new7 = ('[[PERSON_101]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_101]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_101]]' + [char]8220)
$found7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
if (-not $found7) { $failures = $failures + 1 }
Write-Output "Replace 7: $found7"

$old8 = ('[[PERSON_104]] ' + [char]8211 + ' ' + [char]8222 + 'o [[PERSON_104]]' + [char]8220 + ', ' + [char]8222 + 's [[PERSON_104]]' + [char]8220)
$new8 = ('[[PERSON_102]] ' + [char]8211 + ' ' + [char]8222 + 'o [[PERSON_102]]' + [char]8220 + ', ' + [char]8222 + 's [[PERSON_102]]' + [char]8220)
$found8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)
if (-not $found8) { $failures = $failures + 1 }
Write-Output "Replace 8: $found8"

$old9 = ('[[PERSON_105]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_106]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_107]]' + [char]8220)
$new9 = ('[[PERSON_103]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_104]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_104]]' + [char]8220)
$found9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
if (-not $found9) { $failures = $failures + 1 }
Write-Output "Replace 9: $found9"

$old10 = ('[[PERSON_108]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_108]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_108]]' + [char]8220)
$new10 = ('[[PERSON_105]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_105]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_105]]' + [char]8220)
$found10 = $d.Content.Find.Execute($old10, $true, $false, $false, $false, $false, $true, 1, $false, $new10, 2)
if (-not $found10) { $failures = $failures + 1 }
Write-Output "Replace 10: $found10"

$old11 = ('[[PERSON_109]] ' + [char]8211 + ' ' + [char]8222 + 'o [[PERSON_109]]' + [char]8220 + ', ' + [char]8222 + 's [[PERSON_109]]' + [char]8220)
$new11 = ('[[PERSON_106]] ' + [char]8211 + ' ' + [char]8222 + 'o [[PERSON_106]]' + [char]8220 + ', ' + [char]8222 + 's [[PERSON_106]]' + [char]8220)
$found11 = $d.Content.Find.Execute($old11, $true, $false, $false, $false, $false, $true, 1, $false, $new11, 2)
if (-not $found11) { $failures = $failures + 1 }
Write-Output "Replace 11: $found11"

$old12 = ('[[PERSON_110]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_111]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_111]]' + [char]8220)
$new12 = ('[[PERSON_107]] ' + [char]8211 + ' ' + [char]8222 + 's [[PERSON_108]]' + [char]8220 + ', ' + [char]8222 + 'o [[PERSON_108]]' + [char]8220)
$found12 = $d.Content.Find.Execute($old12, $true, $false, $false, $false, $false, $true, 1, $false, $new12, 2)
if (-not $found12) { $failures = $failures + 1 }
Write-Output "Replace 12: $found12"

$old13 = ('[[PERSON_112]] ' + [char]8211 + ' ' + [char]8222 + 'o [[PERSON_112]]' + [char]8220 + ', ' + [char]8222 + 's [[PERSON_112]]' + [char]8220)
$new13 = ('[[PERSON_109]] ' + [char]8211 + ' ' + [char]8222 + 'o [[PERSON_109]]' + [char]8220 + ', ' + [char]8222 + 's [[PERSON_109]]' + [char]8220)
$found13 = $d.Content.Find.Execute($old13, $true, $false, $false, $false, $false, $true, 1, $false, $new13, 2)
if (-not $found13) { $failures = $failures + 1 }
Write-Output "Replace 13: $found13"

Write-Output "Total failures: $failures"